# ACINS-1679 Importer spec completed and bugs fixed.
#
# The course code in row 5 ("AFAM 102.001") was corrected to the "C"
# section ("AFAM 102C.001"), and the active selection was left on that
# row (A5) instead of the header row (A2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the course/section label in A5.
$ws.Range("A5").Value = "AFAM 102C.001"

# Leave the cell selection on A5 (matches the saved workbook view).
$ws.Range("A5").Select() | Out-Null
